# Doing Updates for Financials
# A new reporting period (FY ending 2018-12-31) was added as the first data
# column. This is modeled as: insert a new blank column before column D
# (shifting existing D:K data to E:L), copy the number formats from the
# column that used to be D (now E) into the freshly inserted D column, and
# finally populate the new D column with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column at D; everything from D:K shifts right to E:L.
$ws.Columns("D:D").Insert()

# 2) The newly inserted column D is blank/default-styled. Copy the cell
#    formatting (number format, font, alignment) from column E (which holds
#    what used to be column D) so the new column matches the existing
#    date-row / number-row styles.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# 3) Fill in the new (most recent) period's values in column D.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 291700
$ws.Range("D15").Value = -1300
$ws.Range("D17").Value = 24600
$ws.Range("D18").Value = 267100
$ws.Range("D20").Value = -88900
$ws.Range("D21").Value = 190700
$ws.Range("D23").Value = 178200
$ws.Range("D24").Value = 28200
$ws.Range("D26").Value = 150000
$ws.Range("D27").Value = 150000
$ws.Range("D29").Value = 700
$ws.Range("D32").Value = 88900
$ws.Range("D33").Value = 150600
$ws.Range("D35").Value = 150600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 207800
$ws.Range("D42").Value = 42300
$ws.Range("D48").Value = 133400
$ws.Range("D49").Value = 174700
$ws.Range("D54").Value = 7731900
$ws.Range("D59").Value = 14200
$ws.Range("D66").Value = 6678600
$ws.Range("D72").Value = 606700
$ws.Range("D76").Value = 1053300
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 150600
$ws.Range("D83").Value = 12500
$ws.Range("D89").Value = 189000
$ws.Range("D91").Value = -17600
$ws.Range("D94").Value = -275900
$ws.Range("D96").Value = -53900
$ws.Range("D100").Value = -36800
$ws.Range("D102").Value = -123700
